# Append 4 new match rows (75-78) to the end of the HNL 2023-2024 sheet,
# mirroring the existing layout/format of row 74 (the previous last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (74) down into the
# four new rows so the styled columns (A = bold/bordered index, E = date)
# keep using the same cell styles as the rest of the sheet.
$ws.Range("A74:V74").Copy()
$ws.Range("A75:V75").PasteSpecial(-4122)
$ws.Range("A74:V74").Copy()
$ws.Range("A76:V76").PasteSpecial(-4122)
$ws.Range("A74:V74").Copy()
$ws.Range("A77:V77").PasteSpecial(-4122)
$ws.Range("A74:V74").Copy()
$ws.Range("A78:V78").PasteSpecial(-4122)

$rows = @(
    @{
        r=75; idx=74; f="Varazdin"; g=0; h="Rijeka"; i=2; e=45255.625
        j=5.17;  k="23/11/2023 11:42"; l=5.47; m="25/11/2023 14:59"
        n=3.91;  o="23/11/2023 11:42"; p=4.01; q="25/11/2023 14:59"
        r2=1.59; s="23/11/2023 11:42"; t=1.62; u="25/11/2023 14:59"
        v="https://www.betexplorer.com/football/croatia/hnl/varazdin-rijeka/CSBdZheb/"
    },
    @{
        r=76; idx=75; f="D. Zagreb"; g=2; h="Osijek"; i=1; e=45255.72916666666
        j=1.38;  k="25/11/2023 06:42"; l=1.36; m="25/11/2023 17:22"
        n=4.72;  o="25/11/2023 06:42"; p=5.06; q="25/11/2023 17:26"
        r2=7.99; s="25/11/2023 06:42"; t=8.529999999999999; u="25/11/2023 17:26"
        v="https://www.betexplorer.com/football/croatia/hnl/din-zagreb-osijek/jJA0YCA4/"
    },
    @{
        r=77; idx=76; f="Istra 1961"; g=0; h="Hajduk Split"; i=2; e=45256.625
        j=5.62;  k="22/11/2023 19:42"; l=5.76; m="26/11/2023 14:59"
        n=3.84;  o="22/11/2023 19:42"; p=3.6;  q="26/11/2023 14:59"
        r2=1.61; s="22/11/2023 19:42"; t=1.68; u="26/11/2023 14:53"
        v="https://www.betexplorer.com/football/croatia/hnl/istra-1961-hajduk-split/t83DVAfN/"
    },
    @{
        r=78; idx=77; f="Gorica"; g=3; h="Rudes"; i=0; e=45256.71527777778
        j=1.55;  k="20/11/2023 11:42"; l=1.45; m="26/11/2023 17:03"
        n=3.97;  o="20/11/2023 11:42"; p=4.37; q="26/11/2023 17:07"
        r2=6.15; s="20/11/2023 11:42"; t=7.61; u="26/11/2023 17:08"
        v="https://www.betexplorer.com/football/croatia/hnl/hnk-gorica-rudes/no49WjuH/"
    }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value  = $row.idx      # A - Indice
    $ws.Cells.Item($r, 2).Value  = "croatia"     # B - pais
    $ws.Cells.Item($r, 3).Value  = "hnl"         # C - torneio
    $ws.Cells.Item($r, 4).Value  = "2023-2024"   # D - temporada
    $ws.Cells.Item($r, 5).Value  = $row.e        # E - data_partida
    $ws.Cells.Item($r, 6).Value  = $row.f        # F - home
    $ws.Cells.Item($r, 7).Value  = $row.g        # G - home_ft_gols
    $ws.Cells.Item($r, 8).Value  = $row.h        # H - away
    $ws.Cells.Item($r, 9).Value  = $row.i        # I - away_ft_gols
    $ws.Cells.Item($r, 10).Value = $row.j        # J - home_opening_odds
    $ws.Cells.Item($r, 11).Value = $row.k        # K - home_opening_data_hora
    $ws.Cells.Item($r, 12).Value = $row.l        # L - home_closing_odds
    $ws.Cells.Item($r, 13).Value = $row.m        # M - home_closing_data_hora
    $ws.Cells.Item($r, 14).Value = $row.n        # N - draw_opening_odds
    $ws.Cells.Item($r, 15).Value = $row.o        # O - draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value = $row.p        # P - draw_closing_odds
    $ws.Cells.Item($r, 17).Value = $row.q        # Q - draw_closing_data_hora
    $ws.Cells.Item($r, 18).Value = $row.r2       # R - away_opening_odds
    $ws.Cells.Item($r, 19).Value = $row.s        # S - away_opening_data_hora
    $ws.Cells.Item($r, 20).Value = $row.t        # T - away_closing_odds
    $ws.Cells.Item($r, 21).Value = $row.u        # U - away_closing_data_hora
    $ws.Cells.Item($r, 22).Value = $row.v        # V - url_partida
}

Write-Host "Added rows 75-78"
